$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (2015年) ---
$ws.Range("B4").Value = 6.3283883897
$ws.Range("C4").Value = 32.3374602552
$ws.Range("D4").Value = 7.1504456188
$ws.Range("E4").Value = 7.1138842404
$ws.Range("F4").Value = 159.5097085919
$ws.Range("G4").Value = 23.090751087
$ws.Range("H4").Value = 90.33612469169999
$ws.Range("I4").Value = 8.29643443
$ws.Range("J4").Value = 10.0747204595
$ws.Range("K4").Value = 1.3232692366

# --- Row 5 (2016年) ---
$ws.Range("B5").Value = 6.629517
$ws.Range("C5").Value = 36.813159
$ws.Range("D5").Value = 7.488134
$ws.Range("E5").Value = 7.908716
$ws.Range("F5").Value = 157.236327
$ws.Range("G5").Value = 22.708032
$ws.Range("H5").Value = 91.45995000000001
$ws.Range("I5").Value = 8.484249
$ws.Range("J5").Value = 10.187859
$ws.Range("K5").Value = 1.357324

# --- Row 6 (2017年) ---
$ws.Range("B6").Value = 6.8971104791
$ws.Range("C6").Value = 38.3694385819
$ws.Range("D6").Value = 7.4221984723
$ws.Range("E6").Value = 7.8709119433
$ws.Range("F6").Value = 154.6364888178
$ws.Range("G6").Value = 23.6061410077
$ws.Range("H6").Value = 90.19042291780001
$ws.Range("I6").Value = 8.9368747739
$ws.Range("J6").Value = 10.1103700675
$ws.Range("K6").Value = 1.3678214456

# --- Row 7 (2018年) ---
$ws.Range("B7").Value = 6.9050770176
$ws.Range("C7").Value = 39.875313425
$ws.Range("D7").Value = 7.7636922751
$ws.Range("E7").Value = 7.9858981825
$ws.Range("F7").Value = 148.5335913812
$ws.Range("G7").Value = 27.4670007105
$ws.Range("H7").Value = 87.48129173220001
$ws.Range("I7").Value = 8.4032065752
$ws.Range("J7").Value = 9.9145707364
$ws.Range("K7").Value = 1.3421019716

# --- Row 8 (2019年) ---
$ws.Range("B8").Value = 7.2625067211
$ws.Range("C8").Value = 43.2869775057
$ws.Range("D8").Value = 9.5932565294
$ws.Range("E8").Value = 10.0093260245
$ws.Range("F8").Value = 154.7854278271
$ws.Range("G8").Value = 24.6534143662
$ws.Range("H8").Value = 89.46981795240001
$ws.Range("I8").Value = 9.6070900975
$ws.Range("J8").Value = 9.8233185702
$ws.Range("K8").Value = 1.426710688

# --- Row 9 (2020年) ---
$ws.Range("B9").Value = 7.35011374457948
$ws.Range("C9").Value = 43.8325458438849
$ws.Range("D9").Value = 10.3147448792403
$ws.Range("E9").Value = 12.4335560392008
$ws.Range("I9").Value = 11.7931860222977
$ws.Range("K9").Value = 1.43827353610753

# --- Row 10 (new, 2021年) ---
$ws.Range("A10").Value = "2021年"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("B10").Value = 9.300000000000001
$ws.Range("C10").Value = 52.4
$ws.Range("D10").Value = 10.9
$ws.Range("E10").Value = 12.4
$ws.Range("F10").Value = 170.8
$ws.Range("G10").Value = 30.9
$ws.Range("H10").Value = 107
$ws.Range("I10").Value = 13
$ws.Range("J10").Value = 11.7
$ws.Range("K10").Value = 1.5
